# The document carried a handful of orphaned SharePoint "document library"
# Custom XML Parts (the ct:contentTypeSchema item, the FormTemplates item,
# the empty p:properties item, and their itemProps datastore companions).
# These aren't surfaced anywhere in the document UI/content - they're pure
# package-level metadata left over from a SharePoint export - so clean them
# up the same way Word's own automation model does: walk the
# CustomXMLParts collection back-to-front (so indices stay valid as items
# are removed) and delete every part.

$d = $word.ActiveDocument

$customXmlParts = $d.CustomXMLParts
for ($i = $customXmlParts.Count; $i -ge 1; $i--) {
    $customXmlParts.Item($i).Delete()
}
